$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1022.6111
$ws.Range("I2").Value = 958.0714
$ws.Range("K2").Value = 958.0714
$ws.Range("M2").Value = -845.0714
$ws.Range("H5").Value = 41.2
$ws.Range("I5").Value = 49
$ws.Range("K5").Value = 49
$ws.Range("M5").Value = 66
$ws.Range("H12").Value = 585.45
$ws.Range("I12").Value = 393.07693
$ws.Range("K12").Value = 393.07693
$ws.Range("M12").Value = -223.07693
$ws.Range("H62").Value = 51039.184
$ws.Range("I62").Value = 81322.30499999999
$ws.Range("K62").Value = 81322.30499999999
$ws.Range("M62").Value = -80698.30499999999
$ws.Range("H65").Value = 51039.184
$ws.Range("I65").Value = 81322.30499999999
$ws.Range("K65").Value = 406611.525
$ws.Range("M65").Value = -403491.525
$ws.Range("H86").Value = 4144.8
$ws.Range("I86").Value = 3381.5454
$ws.Range("K86").Value = 3381.5454
$ws.Range("M86").Value = -2258.5454
$ws.Range("H89").Value = 4144.8
$ws.Range("I89").Value = 3381.5454
$ws.Range("K89").Value = 16907.727
$ws.Range("M89").Value = -11291.727
$ws.Range("H98").Value = 26394.75
$ws.Range("I98").Value = 24451.643
$ws.Range("K98").Value = 24451.643
$ws.Range("M98").Value = -22953.643
$ws.Range("H99").Value = 786.6667
$ws.Range("I99").Value = 190.75
$ws.Range("K99").Value = 572.25
$ws.Range("M99").Value = 925.75
$ws.Range("H106").Value = 15900.182
$ws.Range("I106").Value = 16238.375
$ws.Range("K106").Value = 16238.375
$ws.Range("M106").Value = -15607.375
$ws.Range("H107").Value = 167949
$ws.Range("I107").Value = 167949
$ws.Range("K107").Value = 167949
$ws.Range("M107").Value = -166029
$ws.Range("H122").Value = 26394.75
$ws.Range("I122").Value = 24451.643
$ws.Range("K122").Value = 73354.929
$ws.Range("M122").Value = -70904.929
$ws.Range("H137").Value = 31066.236
$ws.Range("I137").Value = 43185.69
$ws.Range("J137").Value = 1124.0588
$ws.Range("K137").Value = 129557.07
$ws.Range("L137").Value = 3372.1764
$ws.Range("M137").Value = -127007.07
$ws.Range("N137").Value = -8472.1764

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 46887.523
$ws.Range("I45").Value = 61077
$ws.Range("J45").Value = 6684
$ws.Range("K45").Value = 61077
$ws.Range("L45").Value = 6684
$ws.Range("M45").Value = -60700
$ws.Range("N45").Value = -7438
$ws.Range("H74").Value = 45884.625
$ws.Range("I74").Value = 8106.39
$ws.Range("K74").Value = 8106.39
$ws.Range("M74").Value = -7232.39
$ws.Range("H77").Value = 45884.625
$ws.Range("I77").Value = 8106.39
$ws.Range("K77").Value = 40531.95
$ws.Range("M77").Value = -36163.95
$ws.Range("H88").Value = 2001.625
$ws.Range("I88").Value = 2000
$ws.Range("J88").Value = 2001.8572
$ws.Range("K88").Value = 2000
$ws.Range("L88").Value = 2001.8572
$ws.Range("M88").Value = -1594
$ws.Range("N88").Value = -2813.8572
$ws.Range("H91").Value = 2001.625
$ws.Range("I91").Value = 2000
$ws.Range("J91").Value = 2001.8572
$ws.Range("K91").Value = 2000
$ws.Range("L91").Value = 2001.8572
$ws.Range("M91").Value = -596
$ws.Range("N91").Value = -4809.8572
$ws.Range("H122").Value = 3260.1428
$ws.Range("I122").Value = 2343.923
$ws.Range("J122").Value = 4749
$ws.Range("K122").Value = 7031.768999999999
$ws.Range("L122").Value = 14247
$ws.Range("M122").Value = -4581.768999999999
$ws.Range("N122").Value = -19147
$ws.Range("H127").Value = 0
$ws.Range("J127").Value = 0
$ws.Range("L127").Value = 0
$ws.Range("H132").Value = 7630.9546
$ws.Range("I132").Value = 7230.8423
$ws.Range("K132").Value = 21692.5269
$ws.Range("M132").Value = -19162.5269

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 26343
$ws.Range("I86").Value = 34882.6
$ws.Range("J86").Value = 4994
$ws.Range("K86").Value = 34882.6
$ws.Range("L86").Value = 4994
$ws.Range("M86").Value = -33759.6
$ws.Range("N86").Value = -7240
$ws.Range("H89").Value = 26343
$ws.Range("I89").Value = 34882.6
$ws.Range("J89").Value = 4994
$ws.Range("K89").Value = 174413
$ws.Range("L89").Value = 24970
$ws.Range("M89").Value = -168797
$ws.Range("N89").Value = -36202
$ws.Range("H134").Value = 7943.857
$ws.Range("I134").Value = 7683.9414
$ws.Range("J134").Value = 8345.546
$ws.Range("K134").Value = 23051.8242
$ws.Range("L134").Value = 25036.638
$ws.Range("M134").Value = -20516.8242
$ws.Range("N134").Value = -30106.638

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 7663.1113
$ws.Range("I6").Value = 7195.8
$ws.Range("K6").Value = 7195.8
$ws.Range("M6").Value = -7082.8
$ws.Range("H9").Value = 75318.53999999999
$ws.Range("J9").Value = 75318.53999999999
$ws.Range("L9").Value = 75318.53999999999
$ws.Range("N9").Value = -75654.53999999999
$ws.Range("H22").Value = 1183.0476
$ws.Range("I22").Value = 703.5
$ws.Range("J22").Value = 2142.1428
$ws.Range("K22").Value = 703.5
$ws.Range("L22").Value = 2142.1428
$ws.Range("M22").Value = -353.5
$ws.Range("N22").Value = -2842.1428
$ws.Range("H31").Value = 25002.936
$ws.Range("I31").Value = 11270.637
$ws.Range("J31").Value = 29198.916
$ws.Range("K31").Value = 11270.637
$ws.Range("L31").Value = 29198.916
$ws.Range("M31").Value = -10975.637
$ws.Range("N31").Value = -29788.916
$ws.Range("H34").Value = 25002.936
$ws.Range("I34").Value = 11270.637
$ws.Range("J34").Value = 29198.916
$ws.Range("K34").Value = 11270.637
$ws.Range("L34").Value = 29198.916
$ws.Range("M34").Value = -11068.637
$ws.Range("N34").Value = -29602.916
$ws.Range("H50").Value = 6532.278
$ws.Range("J50").Value = 6532.278
$ws.Range("L50").Value = 6532.278
$ws.Range("N50").Value = -7782.278
$ws.Range("H132").Value = 57977.168
$ws.Range("I132").Value = 57977.168
$ws.Range("K132").Value = 173931.504
$ws.Range("M132").Value = -171401.504

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 41533.727
$ws.Range("I12").Value = 126991.43
$ws.Range("K12").Value = 380974.29
$ws.Range("M12").Value = -380801.29
$ws.Range("H29").Value = 71.09090999999999
$ws.Range("I29").Value = 20.75
$ws.Range("J29").Value = 99.85714
$ws.Range("K29").Value = 62.25
$ws.Range("L29").Value = 299.57142
$ws.Range("M29").Value = 214.75
$ws.Range("N29").Value = -853.57142
$ws.Range("H68").Value = 2292.7144
$ws.Range("I68").Value = 2592.7273
$ws.Range("J68").Value = 1192.6666
$ws.Range("K68").Value = 7778.1819
$ws.Range("L68").Value = 3577.9998
$ws.Range("M68").Value = -6967.1819
$ws.Range("N68").Value = -5199.9998
$ws.Range("H71").Value = 2292.7144
$ws.Range("I71").Value = 2592.7273
$ws.Range("J71").Value = 1192.6666
$ws.Range("K71").Value = 23334.5457
$ws.Range("L71").Value = 10733.9994
$ws.Range("M71").Value = -19278.5457
$ws.Range("N71").Value = -18845.9994
$ws.Range("H122").Value = 713.2857
$ws.Range("I122").Value = 631.6667
$ws.Range("J122").Value = 774.5
$ws.Range("K122").Value = 5685.0003
$ws.Range("L122").Value = 6970.5
$ws.Range("M122").Value = -3235.0003
$ws.Range("N122").Value = -11870.5

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 11367.454
$ws.Range("I80").Value = 2161.5
$ws.Range("K80").Value = 2161.5
$ws.Range("M80").Value = -1163.5
$ws.Range("H83").Value = 11367.454
$ws.Range("I83").Value = 2161.5
$ws.Range("K83").Value = 10807.5
$ws.Range("M83").Value = -5815.5
$ws.Range("H126").Value = 3808.625
$ws.Range("J126").Value = 4997.857
$ws.Range("L126").Value = 14993.571
$ws.Range("N126").Value = -19933.571
$ws.Range("H136").Value = 56046
$ws.Range("J136").Value = 56046
$ws.Range("L136").Value = 168138
$ws.Range("N136").Value = -173238

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5947.2
$ws.Range("I7").Value = 2118.75
$ws.Range("K7").Value = 2118.75
$ws.Range("M7").Value = -2006.75
$ws.Range("H40").Value = 7017
$ws.Range("I40").Value = 5428.56
$ws.Range("K40").Value = 5428.56
$ws.Range("M40").Value = -5292.56
$ws.Range("H126").Value = 5947.2
$ws.Range("I126").Value = 2118.75
$ws.Range("K126").Value = 6356.25
$ws.Range("M126").Value = -3886.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 3256.9714
$ws.Range("I126").Value = 2653.8518
$ws.Range("K126").Value = 7961.555399999999
$ws.Range("M126").Value = -5491.555399999999
$ws.Range("H136").Value = 5368.39
$ws.Range("I136").Value = 7258
$ws.Range("J136").Value = 2415.875
$ws.Range("K136").Value = 21774
$ws.Range("L136").Value = 7247.625
$ws.Range("M136").Value = -19224
$ws.Range("N136").Value = -12347.625

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("N127").ClearContents()